{"js": "// Bump the Java line numbers embedded in the stack-trace text of the\n// document (template engine moved from 2.0.0 to 2.0.1, shifting the\n// line numbers reported in the exception stack trace).\nconst replacements = [\n  [\"M2DocEvaluator.java:543\", \"M2DocEvaluator.java:555\"],\n  [\"M2DocEvaluator.java:1084\", \"M2DocEvaluator.java:1096\"],\n  [\"M2DocEvaluator.java:1300\", \"M2DocEvaluator.java:1305\"],\n  [\"M2DocEvaluator.java:278\", \"M2DocEvaluator.java:283\"],\n  [\"M2DocEvaluator.java:267\", \"M2DocEvaluator.java:272\"],\n  [\"AbstractTemplatesTestSuite.java:475\", \"AbstractTemplatesTestSuite.java:479\"],\n  [\"AbstractTemplatesTestSuite.java:384\", \"AbstractTemplatesTestSuite.java:388\"],\n];\n\nconst body = context.document.body;\n\nfor (const [searchText, replacementText] of replacements) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replacementText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Bump the Java line numbers embedded in the stack-trace text of the\n# document (template engine moved from 2.0.0 to 2.0.1, shifting the\n# line numbers reported in the exception stack trace).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"M2DocEvaluator.java:543\", \"M2DocEvaluator.java:555\"),\n    @(\"M2DocEvaluator.java:1084\", \"M2DocEvaluator.java:1096\"),\n    @(\"M2DocEvaluator.java:1300\", \"M2DocEvaluator.java:1305\"),\n    @(\"M2DocEvaluator.java:278\", \"M2DocEvaluator.java:283\"),\n    @(\"M2DocEvaluator.java:267\", \"M2DocEvaluator.java:272\"),\n    @(\"AbstractTemplatesTestSuite.java:475\", \"AbstractTemplatesTestSuite.java:479\"),\n    @(\"AbstractTemplatesTestSuite.java:384\", \"AbstractTemplatesTestSuite.java:388\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
